# Generate Report for Handback
# Two source-file GUIDs in this handback report are renamed/re-run:
#   6bc5fd6b-83dd-44ab-a817-8de84405dc24  ->  dabeae51-39e2-4f6f-9c0e-ff304514674e
#   9f75aa24-c862-4956-be99-e0355a2c7a3b  ->  ffff387d067a-3f28-42a4-a454-346c8390a718
# and the handoff/handback xliff names + timestamps refresh to reflect a
# fresh CI run. Update every cell that carries one of these derived values,
# on all three sheets, and keep each hyperlink's display text (not its
# target URL) in sync with its cell text.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "6bc5fd6b-83dd-44ab-a817-8de84405dc24"
$oldGuid2 = "9f75aa24-c862-4956-be99-e0355a2c7a3b"
$newGuid1 = "dabeae51-39e2-4f6f-9c0e-ff304514674e"
$newGuid2 = "ffff387d067a-3f28-42a4-a454-346c8390a718"

# ---- Sheet 1: Overview ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-09-03 03:06:51"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-09-03 03:06:51"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid1.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newGuid2.md"
    }
}

# ---- Sheet 2: zh-cn ---------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlfNew = "$newGuid1.c6d2fed708f9c2fb3a547dfdd9af584f04ae08a2.zh-cn.xlf"

$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = $zhXlfNew
$wsZhCn.Range("H2").Value = "2016-09-03 03:06:46"
$wsZhCn.Range("J2").Value = $zhXlfNew
$wsZhCn.Range("K2").Value = "2016-09-03 03:07:08"

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = $zhXlfNew
$wsZhCn.Range("H3").Value = "2016-09-03 03:06:46"
$wsZhCn.Range("J3").Value = $zhXlfNew
$wsZhCn.Range("K3").Value = "2016-09-03 03:07:08"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newGuid2.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}

# ---- Sheet 3: de-de ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlfNew = "$newGuid1.c6d2fed708f9c2fb3a547dfdd9af584f04ae08a2.de-de.xlf"

$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = $deXlfNew
$wsDeDe.Range("H2").Value = "2016-09-03 03:06:51"
$wsDeDe.Range("J2").Value = $deXlfNew
$wsDeDe.Range("K2").Value = "2016-09-03 03:07:16"

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = $deXlfNew
$wsDeDe.Range("H3").Value = "2016-09-03 03:06:51"
$wsDeDe.Range("J3").Value = $deXlfNew
$wsDeDe.Range("K3").Value = "2016-09-03 03:07:16"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newGuid2.md"
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}
